# Updated cryptos list on Mon Aug 12 02:32:13 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with newly scraped
# values for every coin row, and re-sorts three coin pairs whose rank
# order flipped (Aptos/Monero, Bittensor/Filecoin, WhiteBITCoin/EnergySwap)
# by swapping their Coin (B) / Link (C) cells too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "58.451.82" },
    @{ Cell = "E2"; Value = "  -4.14%  " },
    @{ Cell = "D3"; Value = "2.534.07" },
    @{ Cell = "E3"; Value = "  -3.29%  " },
    @{ Cell = "E4"; Value = "  +0.02%  " },
    @{ Cell = "D5"; Value = "506.40" },
    @{ Cell = "E5"; Value = "  -4.10%  " },
    @{ Cell = "D6"; Value = "143.08" },
    @{ Cell = "E6"; Value = "  -7.72%  " },
    @{ Cell = "E7"; Value = "  +0.07%  " },
    @{ Cell = "E8"; Value = "  -5.39%  " },
    @{ Cell = "D9"; Value = "2.535.98" },
    @{ Cell = "E9"; Value = "  -3.45%  " },
    @{ Cell = "D10"; Value = "6.15" },
    @{ Cell = "E10"; Value = "  -8.03%  " },
    @{ Cell = "E11"; Value = "  -5.79%  " },
    @{ Cell = "D12"; Value = "0.330" },
    @{ Cell = "E12"; Value = "  -4.85%  " },
    @{ Cell = "E13"; Value = "  -0.61%  " },
    @{ Cell = "D14"; Value = "2.979.78" },
    @{ Cell = "E14"; Value = "  -3.20%  " },
    @{ Cell = "D15"; Value = "58.433.21" },
    @{ Cell = "E15"; Value = "  -4.16%  " },
    @{ Cell = "D16"; Value = "20.63" },
    @{ Cell = "E16"; Value = "  -5.05%  " },
    @{ Cell = "E17"; Value = "  -5.50%  " },
    @{ Cell = "D18"; Value = "2.538.28" },
    @{ Cell = "E18"; Value = "  -3.18%  " },
    @{ Cell = "D19"; Value = "4.52" },
    @{ Cell = "E19"; Value = "  -5.35%  " },
    @{ Cell = "D20"; Value = "334.13" },
    @{ Cell = "E20"; Value = "  -5.78%  " },
    @{ Cell = "D21"; Value = "10.05" },
    @{ Cell = "E21"; Value = "  -5.31%  " },
    @{ Cell = "D22"; Value = "0.998" },
    @{ Cell = "E22"; Value = "  -0.22%  " },
    @{ Cell = "D23"; Value = "5.93" },
    @{ Cell = "E23"; Value = "  -4.91%  " },
    @{ Cell = "D24"; Value = "59.96" },
    @{ Cell = "E24"; Value = "  -2.64%  " },
    @{ Cell = "D25"; Value = "0.406" },
    @{ Cell = "E25"; Value = "  -5.29%  " },
    @{ Cell = "D26"; Value = "0.999" },
    @{ Cell = "E26"; Value = "  -0.05%  " },
    @{ Cell = "E27"; Value = "  -5.14%  " },
    @{ Cell = "D28"; Value = "2.652.40" },
    @{ Cell = "E28"; Value = "  -3.04%  " },
    @{ Cell = "D29"; Value = "0.0₃0784" },
    @{ Cell = "E29"; Value = "  -8.93%  " },
    @{ Cell = "E30"; Value = "  -6.82%  " },
    @{ Cell = "D31"; Value = "1.00" },
    @{ Cell = "E31"; Value = "  +0.03%  " },
    @{ Cell = "B32"; Value = "Monero" },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" },
    @{ Cell = "D32"; Value = "149.40" },
    @{ Cell = "E32"; Value = "  -0.76%  " },
    @{ Cell = "B33"; Value = "Aptos" },
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" },
    @{ Cell = "D33"; Value = "5.84" },
    @{ Cell = "E33"; Value = "  -5.57%  " },
    @{ Cell = "D34"; Value = "18.48" },
    @{ Cell = "E34"; Value = "  -4.94%  " },
    @{ Cell = "E35"; Value = "  -5.13%  " },
    @{ Cell = "D36"; Value = "0.932" },
    @{ Cell = "E36"; Value = "  +4.01%  " },
    @{ Cell = "E37"; Value = "  -7.32%  " },
    @{ Cell = "E38"; Value = "  -7.99%  " },
    @{ Cell = "D39"; Value = "36.01" },
    @{ Cell = "E39"; Value = "  -1.15%  " },
    @{ Cell = "D40"; Value = "0.825" },
    @{ Cell = "E40"; Value = "  -12.23%  " },
    @{ Cell = "D41"; Value = "1.40" },
    @{ Cell = "E41"; Value = "  -6.83%  " },
    @{ Cell = "B42"; Value = "Filecoin" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" },
    @{ Cell = "D42"; Value = "3.52" },
    @{ Cell = "E42"; Value = "  -7.28%  " },
    @{ Cell = "B43"; Value = "Bittensor" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao" },
    @{ Cell = "D43"; Value = "281.56" },
    @{ Cell = "E43"; Value = "  -4.40%  " },
    @{ Cell = "D44"; Value = "0.0994" },
    @{ Cell = "E44"; Value = "  -3.42%  " },
    @{ Cell = "E45"; Value = "  -0.01%  " },
    @{ Cell = "D46"; Value = "0.601" },
    @{ Cell = "E46"; Value = "  -5.78%  " },
    @{ Cell = "E47"; Value = "  -4.69%  " },
    @{ Cell = "B48"; Value = "EnergySwap" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" },
    @{ Cell = "D48"; Value = "18.63" },
    @{ Cell = "E48"; Value = "  -5.57%  " },
    @{ Cell = "B49"; Value = "WhiteBITCoin" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt" },
    @{ Cell = "D49"; Value = "10.29" },
    @{ Cell = "E49"; Value = "  -0.60%  " },
    @{ Cell = "E50"; Value = "  -5.09%  " },
    @{ Cell = "D51"; Value = "4.50" },
    @{ Cell = "E51"; Value = "  -7.76%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force plain-text storage so numeric-looking strings (prices like
    # "1.00", "506.40", ids like "0.330") aren't silently reinterpreted
    # as numbers by Excel's type inference - then drop back to the
    # workbook's default ("Normal") style so no stray per-cell number
    # format is left behind on save.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
